# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" sheet (fund-holdings detail) between the
#   existing "2021-Q4" sheet and the "总计" (summary) sheet.
# - Rebuilds the "总计" summary sheet with an extra row for 2022-Q1
#   (placed above the pre-existing 2021-Q4 row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: style a header cell (bold, centered, top-aligned, thin border)
# ---------------------------------------------------------------------
function Format-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------
# Helper: write a value that must stay TEXT (avoids Excel's automatic
# "looks like a number" conversion for things like fund codes "519019"
# or formatted numeric strings like "10.33"), then drop back to the
# default "Normal" style so no stray number-format styling lingers.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# =======================================================================
# 1. Remove the old "总计" sheet - it will be rebuilt from scratch further
#    down so that its internal sheetId continues the counter *after* the
#    new "2022-Q1" sheet, matching 2021-Q4(1) / 2022-Q1(2) / 总计(3).
# =======================================================================
$wb.Worksheets.Item("总计").Delete()

$sheetQ4 = $wb.Worksheets.Item(1)

# =======================================================================
# 2. Build the new "2022-Q1" sheet right after "2021-Q4".
# =======================================================================
$q1 = $wb.Worksheets.Add($null, $sheetQ4)
$q1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q1.Cells.Item(1, $col)
    Set-TextValue $cell $headers[$col - 2]
    Format-HeaderCell $cell
}

$rows = @(
    @(0,"519019","大成景阳领先混合","10.33","92.80","5.32","0.5496",7),
    @(1,"004040","金鹰医疗健康产业股票A","10.91","92.37","4.96","0.5411",9),
    @(2,"011335","银河医药健康混合型证券投资基金","8.94","87.69","5.59","0.4997",7),
    @(3,"004041","金鹰医疗健康产业股票C","7.31","92.37","4.96","0.3626",9),
    @(4,"090020","大成健康产业混合","3.76","91.73","6.04","0.2271",5),
    @(5,"090016","大成消费主题混合","4.23","93.78","5.26","0.2225",7),
    @(6,"012045","大成医药健康股票A","2.87","93.58","7.47","0.2144",3),
    @(7,"002319","大成一带一路灵活配置混合","0.50","89.30","5.96","0.0298",5),
    @(8,"001861","富安达健康人生灵活配置混合","0.61","82.18","3.12","0.0190",9),
    @(9,"012046","大成医药健康股票C","0.25","93.58","7.47","0.0187",3),
    @(10,"011377","创金合信积极成长股票A","0.29","94.90","4.91","0.0142",5),
    @(11,"002789","长盛同享灵活配置混合A","0.49","79.17","2.35","0.0115",8),
    @(12,"014285","鑫元健康产业混合A","0.40","33.76","1.92","0.0077",9),
    @(13,"011378","创金合信积极成长股票C","0.11","94.90","4.91","0.0054",5),
    @(14,"014286","鑫元健康产业混合C","0.27","33.76","1.92","0.0052",9),
    @(15,"002790","长盛同享灵活配置混合C","0.02","79.17","2.35","0.0005",8)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q1.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Format-HeaderCell $aCell

    Set-TextValue $q1.Cells.Item($r, 2) $row[1]
    Set-TextValue $q1.Cells.Item($r, 3) $row[2]
    Set-TextValue $q1.Cells.Item($r, 4) $row[3]
    Set-TextValue $q1.Cells.Item($r, 5) $row[4]
    Set-TextValue $q1.Cells.Item($r, 6) $row[5]
    Set-TextValue $q1.Cells.Item($r, 7) $row[6]

    $q1.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# =======================================================================
# 3. Rebuild the "总计" (summary) sheet right after "2022-Q1", with a new
#    row for 2022-Q1 on top of the pre-existing 2021-Q4 row.
# =======================================================================
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $cell = $total.Cells.Item(1, $col)
    Set-TextValue $cell $totalHeaders[$col - 2]
    Format-HeaderCell $cell
}

$totalRows = @(
    @(0, "2022-Q1", 16, 2.73),
    @(1, "2021-Q4", 30, 8.949999999999999)
)

$r = 2
foreach ($row in $totalRows) {
    $aCell = $total.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    Format-HeaderCell $aCell

    Set-TextValue $total.Cells.Item($r, 2) $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]

    $r++
}

$wb.Worksheets.Item(1).Select()
